$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.454.53'
$ws.Range('E2').Value = '  -4.64%  '
# Row 3
$ws.Range('D3').Value = '3.278.50'
$ws.Range('E3').Value = '  -7.15%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.18'
$ws.Range('E5').Value = '  -4.44%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.42'
$ws.Range('E6').Value = '  -11.70%  '
# Row 7
$ws.Range('E7').Value = '  +0.01%  '
# Row 8
$ws.Range('D8').Value = '3.268.42'
$ws.Range('E8').Value = '  -7.33%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  -10.48%  '
# Row 10
$ws.Range('E10').Value = '  -13.33%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.70'
$ws.Range('E11').Value = '  -7.10%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.512'
$ws.Range('E12').Value = '  -12.65%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.77'
$ws.Range('E13').Value = '  -16.06%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000245'
$ws.Range('E14').Value = '  -10.88%  '
# Row 15
$ws.Range('D15').Value = '3.807.47'
$ws.Range('E15').Value = '  -7.09%  '
# Row 16
$ws.Range('D16').Value = '67.499.44'
$ws.Range('E16').Value = '  -4.72%  '
# Row 17
$ws.Range('D17').Value = '3.280.42'
$ws.Range('E17').Value = '  -7.30%  '
# Row 18
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '537.12'
$ws.Range('E18').Value = '  -11.50%  '
# Row 19
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.27'
$ws.Range('E19').Value = '  -13.80%  '
# Row 20
$ws.Range('E20').Value = '  -6.08%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.14'
$ws.Range('E21').Value = '  -14.31%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.764'
$ws.Range('E22').Value = '  -13.26%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.90'
$ws.Range('E23').Value = '  -12.91%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.71'
$ws.Range('E24').Value = '  -12.45%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.95'
$ws.Range('E25').Value = '  -12.06%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.26'
$ws.Range('E27').Value = '  -11.67%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.16'
$ws.Range('E28').Value = '  -10.59%  '
# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '29.54'
$ws.Range('E29').Value = '  -12.19%  '
# Row 30
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.16'
$ws.Range('E30').Value = '  -15.63%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.67'
$ws.Range('E31').Value = '  -10.70%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -10.83%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.64'
$ws.Range('E33').Value = '  -18.29%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '540.18'
$ws.Range('E34').Value = '  -12.37%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.82'
$ws.Range('E35').Value = '  -14.53%  '
# Row 36
$ws.Range('E36').Value = '  +0.05%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0464'
$ws.Range('E37').Value = '  -6.98%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.52'
$ws.Range('E38').Value = '  -5.98%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0864'
$ws.Range('E39').Value = '  -13.13%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.129'
$ws.Range('E40').Value = '  -9.82%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.14'
$ws.Range('E41').Value = '  -15.70%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.77'
$ws.Range('E42').Value = '  -18.28%  '
# Row 43
$ws.Range('D43').Value = '2.948.25'
$ws.Range('E43').Value = '  -11.94%  '
# Row 44
$ws.Range('E44').Value = '  -12.56%  '
# Row 45
$ws.Range('E45').Value = '  -17.39%  '
# Row 46
$ws.Range('E46').Value = '  -10.95%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.94'
$ws.Range('E47').Value = '  -15.43%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  -18.10%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.78'
# Row 51
$ws.Range('E51').Value = '  -11.86%  '
